$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -12.37039999999999
    21 = -13.2884
    23 = -12.3529
    25 = -11.7357
    53 = -11.16090000000001
    57 = -14.21479999999999
    59 = -12.638
    69 = -11.10509999999999
    79 = -11.55
    83 = -13.53669999999999
    93 = -10.0772
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
